$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.513.34"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.555.70"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").Value = "3.555.35"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +3.32%  "
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.415"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "4.158.16"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "3.556.20"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "66.559.11"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.614"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("D25").Value = "3.695.80"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").Value = "3.549.79"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.154"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("B37").Value = "USDe"
$ws.Range("C37").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "172.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0848"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.887"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.58%  "
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.940"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.52%  "
